# Refresh the legacy GSC export ("HTTPS" / Chart sheet):
#   - drop the oldest day (2025-10-15) from the top of the table
#   - shift every remaining day/value up by one row
#   - append the newest day (2026-01-13) at the bottom with a 0 count
#
# Done with Range Copy / PasteSpecial (values only) so the shifted date
# strings keep being stored as plain text (same as the source cells),
# instead of Excel's "looks like a date -> convert to a date serial"
# auto-conversion that a direct .Value assignment would trigger.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = 91
$xlPasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues

# 1) Shift rows 3..91 up into rows 2..90 (dates in col A, counts in col C;
#    col B is always 0 so it rides along for free).
$source = $ws.Range("A3:C" + $lastRow)
$dest = $ws.Range("A2:C" + ($lastRow - 1))
$source.Copy()
$dest.PasteSpecial($xlPasteValues)

# 2) Stage the brand-new date as a text formula result (so it is never
#    parsed as literal user input / auto-converted to a date serial),
#    then copy *that cell's value* into the new last row so it lands as
#    a normal shared-string text cell, matching the rest of column A.
$helper = $ws.Range("Z1")
$helper.Formula = "=""2026-01-13"""
$helper.Copy()
$ws.Range("A" + $lastRow).PasteSpecial($xlPasteValues)
$helper.ClearContents()

# 3) The newest day has no recorded HTTPS URL count yet.
$ws.Range("B" + $lastRow).Value = 0
$ws.Range("C" + $lastRow).Value = 0

$excel.CutCopyMode = 0
